$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 181.26666
$ws.Range("I39").Value = 173.33333
$ws.Range("J39").Value = 213
$ws.Range("K39").Value = 519.99999
$ws.Range("L39").Value = 639
$ws.Range("M39").Value = -223.99999
$ws.Range("N39").Value = -1231
$ws.Range("H41").Value = 4527.0835
$ws.Range("J41").Value = 6008.3335
$ws.Range("L41").Value = 6008.3335
$ws.Range("N41").Value = -6888.3335
$ws.Range("H43").Value = 1090
$ws.Range("I43").Value = 950
$ws.Range("K43").Value = 950
$ws.Range("M43").Value = -881
$ws.Range("H111").Value = 12927.211
$ws.Range("I111").Value = 914.7778
$ws.Range("J111").Value = 23738.4
$ws.Range("K111").Value = 2744.3334
$ws.Range("L111").Value = 71215.20000000001
$ws.Range("M111").Value = 322.6666
$ws.Range("N111").Value = -77349.20000000001
$ws.Range("H116").Value = 53147480
$ws.Range("I116").Value = 62781148
$ws.Range("K116").Value = 62781148
$ws.Range("M116").Value = -62777706
$ws.Range("H131").Value = 13653.77
$ws.Range("I131").Value = 3633.3333
$ws.Range("J131").Value = 16659.9
$ws.Range("K131").Value = 10899.9999
$ws.Range("L131").Value = 49979.7
$ws.Range("M131").Value = -5859.999899999999
$ws.Range("N131").Value = -60059.7
$ws.Range("H138").Value = 2007.711
$ws.Range("J138").Value = 2915.8
$ws.Range("L138").Value = 8747.400000000001
$ws.Range("N138").Value = -19027.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2048.4
$ws.Range("I45").Value = 1942.6666
$ws.Range("K45").Value = 1942.6666
$ws.Range("M45").Value = -1565.6666
$ws.Range("H63").Value = 133338504
$ws.Range("I63").Value = 250001630
$ws.Range("J63").Value = 40008000
$ws.Range("K63").Value = 250001630
$ws.Range("L63").Value = 40008000
$ws.Range("M63").Value = -250000944
$ws.Range("N63").Value = -40009372
$ws.Range("H66").Value = 133338504
$ws.Range("I66").Value = 250001630
$ws.Range("J66").Value = 40008000
$ws.Range("K66").Value = 1250008150
$ws.Range("L66").Value = 200040000
$ws.Range("M66").Value = -1250004718
$ws.Range("N66").Value = -200046864
$ws.Range("H122").Value = 7409990
$ws.Range("I122").Value = 10103084
$ws.Range("J122").Value = 3983
$ws.Range("K122").Value = 30309252
$ws.Range("L122").Value = 11949
$ws.Range("M122").Value = -30306802
$ws.Range("N122").Value = -16849

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1224
$ws.Range("I94").Value = 867.5
$ws.Range("K94").Value = 867.5
$ws.Range("M94").Value = -416.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2655.6
$ws.Range("I6").Value = 3427
$ws.Range("K6").Value = 3427
$ws.Range("M6").Value = -3314
$ws.Range("H16").Value = 624.8
$ws.Range("I16").Value = 624.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 624.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -337.8
$ws.Range("N16").ClearContents()
$ws.Range("H99").Value = 2527.9565
$ws.Range("I99").Value = 1982.7858
$ws.Range("K99").Value = 1982.7858
$ws.Range("M99").Value = -484.7858000000001
$ws.Range("H113").Value = 624.8
$ws.Range("I113").Value = 624.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 624.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1545.2
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2527.9565
$ws.Range("I126").Value = 1982.7858
$ws.Range("K126").Value = 5948.357400000001
$ws.Range("M126").Value = -3478.357400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 396.83334
$ws.Range("I5").Value = 396.83334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1190.50002
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1078.50002
$ws.Range("N5").ClearContents()
$ws.Range("H29").Value = 335
$ws.Range("I29").Value = 335
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1005
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -728
$ws.Range("N29").ClearContents()
$ws.Range("H135").Value = 396.83334
$ws.Range("I135").Value = 396.83334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3571.50006
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1036.50006
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 2619.0476
$ws.Range("I136").Value = 6333.3335
$ws.Range("K136").Value = 19000.0005
$ws.Range("M136").Value = -13900.0005
$ws.Range("H137").Value = 675
$ws.Range("I137").Value = 675
$ws.Range("K137").Value = 2025
$ws.Range("M137").Value = 3075
$ws.Range("H138").Value = 5555859
$ws.Range("I138").Value = 8572809
$ws.Range("J138").Value = 276197
$ws.Range("K138").Value = 25718427
$ws.Range("L138").Value = 828591
$ws.Range("M138").Value = -25713287
$ws.Range("N138").Value = -838871
$ws.Range("H141").Value = 2787.6667
$ws.Range("I141").Value = 2787.6667
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8363.000100000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3183.000100000001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2924.0908
$ws.Range("I80").Value = 2763.4285
$ws.Range("K80").Value = 2763.4285
$ws.Range("M80").Value = -1765.4285
$ws.Range("H83").Value = 2924.0908
$ws.Range("I83").Value = 2763.4285
$ws.Range("K83").Value = 13817.1425
$ws.Range("M83").Value = -8825.1425
$ws.Range("H102").Value = 2693.9412
$ws.Range("I102").Value = 1379.7
$ws.Range("J102").Value = 4571.4287
$ws.Range("K102").Value = 1379.7
$ws.Range("L102").Value = 4571.4287
$ws.Range("M102").Value = 242.3
$ws.Range("N102").Value = -7815.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4785.7144
$ws.Range("I40").Value = 3750
$ws.Range("K40").Value = 3750
$ws.Range("M40").Value = -3614
$ws.Range("H55").Value = 553.125
$ws.Range("I55").Value = 544.8461
$ws.Range("J55").Value = 589
$ws.Range("K55").Value = 544.8461
$ws.Range("L55").Value = 589
$ws.Range("M55").Value = -371.8461
$ws.Range("N55").Value = -935
$ws.Range("H132").Value = 6571.952
$ws.Range("I132").Value = 3610
$ws.Range("J132").Value = 9264.637000000001
$ws.Range("K132").Value = 10830
$ws.Range("L132").Value = 27793.911
$ws.Range("M132").Value = -8300
$ws.Range("N132").Value = -32853.911
$ws.Range("H134").Value = 38657.25
$ws.Range("J134").Value = 38657.25
$ws.Range("L134").Value = 38657.25
$ws.Range("N134").Value = -48797.25
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 707.4286
$ws.Range("I107").Value = 707.4286
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2122.2858
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -202.2857999999997
$ws.Range("N107").ClearContents()
$ws.Range("H126").Value = 2347.3333
$ws.Range("I126").Value = 2153.25
$ws.Range("K126").Value = 6459.75
$ws.Range("M126").Value = -3989.75

